$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: new "s23" target (was previously the "Odroid-n2" row before insertion)
$ws.Range("A11").Value = "s23"
$ws.Range("B11").Value = "Galaxy_S22"
$ws.Range("C11").Value = "tflite"
$ws.Range("D11").Value = "android"
$ws.Range("E11").Value = "arm"
$ws.Range("F11").Value = "cpu"
$ws.Range("G11").Value = "8G"
$ws.Range("H11").Value = "default_target.png"

# Row 12: Odroid-n2 (shifted down from row 11, same data as before)
$ws.Range("A12").Value = "Odroid-n2"
$ws.Range("B12").Value = "Odroid_N2"
$ws.Range("C12").Value = "acl"
$ws.Range("D12").Value = "ubuntu"
$ws.Range("E12").Value = "arm"
$ws.Range("F12").Value = "opencl"
$ws.Range("G12").Value = "2G"
$ws.Range("H12").Value = "odroid-n2.jpg"

# Row 13: Odroid M1 (new)
$ws.Range("A13").Value = "Odroid M1"
$ws.Range("B13").Value = "Odroid_M1"
$ws.Range("C13").Value = "rknn"
$ws.Range("D13").Value = "ubuntu"
$ws.Range("E13").Value = "arm"
$ws.Range("F13").Value = "opencl"
$ws.Range("G13").Value = "2G"
$ws.Range("H13").Value = "default_target.png"

# Row 14: Rasberry Pi5 (new)
$ws.Range("A14").Value = "Rasberry Pi5"
$ws.Range("B14").Value = "Rasberry_Pi5"
$ws.Range("C14").Value = "tflite"
$ws.Range("D14").Value = "ubuntu"
$ws.Range("E14").Value = "arm"
$ws.Range("F14").Value = "opencl"
$ws.Range("G14").Value = "2G"
$ws.Range("H14").Value = "default_target.png"

# Row 15: Comma 3X (new)
$ws.Range("A15").Value = "Comma 3X"
$ws.Range("B15").Value = "Comma_3X"
$ws.Range("C15").Value = "pytorch"
$ws.Range("D15").Value = "ubuntu"
$ws.Range("E15").Value = "arm"
$ws.Range("F15").Value = "opencl"
$ws.Range("G15").Value = "2G"
$ws.Range("H15").Value = "default_target.png"

# Row 16: KT cloud (new)
$ws.Range("A16").Value = "KT cloud"
$ws.Range("B16").Value = "KT_cloud"
$ws.Range("C16").Value = "pytorch"
$ws.Range("D16").Value = "ubuntu"
$ws.Range("E16").Value = "arm"
$ws.Range("F16").Value = "opencl"
$ws.Range("G16").Value = "2G"
$ws.Range("H16").Value = "default_target.png"

# Row 17: Amazon Web Services (new)
$ws.Range("A17").Value = "Amazon Web Services"
$ws.Range("B17").Value = "AWS"
$ws.Range("C17").Value = "pytorch"
$ws.Range("D17").Value = "ubuntu"
$ws.Range("E17").Value = "x86"
$ws.Range("F17").Value = "cpu"
$ws.Range("G17").Value = "2G"
$ws.Range("H17").Value = "default_target.png"

# Row 18: Google Cloud Platform (new)
$ws.Range("A18").Value = "Google Cloud Platform"
$ws.Range("B18").Value = "GCP"
$ws.Range("C18").Value = "pytorch"
$ws.Range("D18").Value = "ubuntu"
$ws.Range("E18").Value = "x86"
$ws.Range("F18").Value = "cpu"
$ws.Range("G18").Value = "2G"
$ws.Range("H18").Value = "default_target.png"

# Row 18 is slightly shorter than the default row height in the target workbook
$ws.Rows.Item(18).RowHeight = 12

# Column A widened to fit the new, longer target names (e.g. "Amazon Web Services")
$ws.Columns.Item(1).ColumnWidth = 29.28

# Move the selection to match the author's final cursor position
$ws.Range("G26").Select()
